# Generate Report for Handback
# Refresh the handoff/handback timestamp columns that the report regenerates
# each run. These cells hold plain text (formatted to look like dates), so
# assign strings to avoid Excel coercing them into date serials.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: "Latest HO Xliff Generate Date" for the first file
# mirrors the de-de "Correspond Handoff Datetime" for that same file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 15:20:53"

# "zh-cn" sheet: handoff / handback datetimes for the first file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-05 15:20:49"
$wsZhCn.Range("K2").Value = "2016-09-05 15:21:13"

# "de-de" sheet: handoff / handback datetimes for the first file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-05 15:20:53"
$wsDeDe.Range("K2").Value = "2016-09-05 15:21:21"
